$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 208.75
$ws.Range("K33").Value = 208.75
$ws.Range("I33").Value = 208.75
$ws.Range("M33").Value = 20.25
$ws.Range("H70").Value = 1814.6666
$ws.Range("L70").Value = 5250
$ws.Range("J70").Value = 1750
$ws.Range("N70").Value = -5790
$ws.Range("H73").Value = 1814.6666
$ws.Range("L73").Value = 5250
$ws.Range("J73").Value = 1750
$ws.Range("N73").Value = -7122
$ws.Range("H127").Value = 1132
$ws.Range("K127").Value = 1200
$ws.Range("I127").Value = 400
$ws.Range("M127").Value = 3760
$ws.Range("H135").Value = 654
$ws.Range("K135").Value = 5886
$ws.Range("I135").Value = 654
$ws.Range("M135").Value = -3351
$ws.Range("H137").Value = 3375.9285
$ws.Range("L137").Value = 10462.9095
$ws.Range("J137").Value = 3487.6365
$ws.Range("N137").Value = -15562.9095
$ws.Range("H138").Value = 13328.333
$ws.Range("L138").Value = 45268.857
$ws.Range("J138").Value = 15089.619
$ws.Range("N138").Value = -55548.857
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8695.25
$ws.Range("K32").Value = 7942.3687
$ws.Range("I32").Value = 7942.3687
$ws.Range("M32").Value = -7655.3687
$ws.Range("H45").Value = 21053.363
$ws.Range("K45").Value = 6512.4287
$ws.Range("L45").Value = 46500
$ws.Range("I45").Value = 6512.4287
$ws.Range("J45").Value = 46500
$ws.Range("M45").Value = -6135.4287
$ws.Range("N45").Value = -47254
$ws.Range("H63").Value = 10496.667
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 10496.667
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 10496.667
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -11868.667
$ws.Range("H66").Value = 10496.667
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 52483.335
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 10496.667
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -59347.335
$ws.Range("H74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H97").Value = 1322.1818
$ws.Range("K97").Value = 1054.5
$ws.Range("I97").Value = 1054.5
$ws.Range("M97").Value = -558.5
$ws.Range("H119").Value = 106666.336
$ws.Range("L119").Value = 106666.336
$ws.Range("J119").Value = 106666.336
$ws.Range("N119").Value = -116342.336
$ws.Range("H132").Value = 3499.5833
$ws.Range("K132").Value = 8999.25
$ws.Range("I132").Value = 2999.75
$ws.Range("M132").Value = -6469.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 24864.666
$ws.Range("L95").Value = 24864.666
$ws.Range("J95").Value = 24864.666
$ws.Range("N95").Value = -30356.666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2313.8235
$ws.Range("K31").Value = 1894.4546
$ws.Range("L31").Value = 3082.6667
$ws.Range("I31").Value = 1894.4546
$ws.Range("J31").Value = 3082.6667
$ws.Range("M31").Value = -1599.4546
$ws.Range("N31").Value = -3672.6667
$ws.Range("H34").Value = 2313.8235
$ws.Range("K34").Value = 1894.4546
$ws.Range("L34").Value = 3082.6667
$ws.Range("I34").Value = 1894.4546
$ws.Range("J34").Value = 3082.6667
$ws.Range("M34").Value = -1692.4546
$ws.Range("N34").Value = -3486.6667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 299.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1797
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 299.5
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2023
$ws.Range("H7").Value = 107.1
$ws.Range("K7").Value = 337.89474
$ws.Range("L7").Value = 6
$ws.Range("I7").Value = 112.63158
$ws.Range("J7").Value = 2
$ws.Range("M7").Value = -225.89474
$ws.Range("N7").Value = -230
$ws.Range("H11").Value = 1585027.5
$ws.Range("K11").Value = 5364802.800000001
$ws.Range("L11").Value = 1300000.02
$ws.Range("I11").Value = 1788267.6
$ws.Range("J11").Value = 433333.34
$ws.Range("M11").Value = -5364662.800000001
$ws.Range("N11").Value = -1300280.02
$ws.Range("H17").Value = 3000
$ws.Range("L17").Value = 9000
$ws.Range("J17").Value = 3000
$ws.Range("N17").Value = -9338
$ws.Range("H38").Value = 105
$ws.Range("K38").Value = 315
$ws.Range("I38").Value = 105
$ws.Range("M38").Value = 32
$ws.Range("H114").Value = 1334.125
$ws.Range("L114").Value = 3687
$ws.Range("J114").Value = 1229
$ws.Range("N114").Value = -10195
$ws.Range("H131").Value = 1677.5714
$ws.Range("L131").Value = 5142
$ws.Range("J131").Value = 1714
$ws.Range("N131").Value = -15222
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H45").Value = 52000
$ws.Range("L45").Value = 79000
$ws.Range("J45").Value = 79000
$ws.Range("N45").Value = -80118
$ws.Range("H97").Value = 6333
$ws.Range("K97").Value = 3999.5
$ws.Range("I97").Value = 3999.5
$ws.Range("M97").Value = -3503.5
$ws.Range("H123").Value = 120000
$ws.Range("L123").Value = 120000
$ws.Range("J123").Value = 120000
$ws.Range("N123").Value = -124900
$ws.Range("H126").Value = 3868.6
$ws.Range("K126").Value = 11509.5
$ws.Range("I126").Value = 3836.5
$ws.Range("M126").Value = -9039.5
$ws.Range("H132").Value = 4664.273
$ws.Range("K132").Value = 11918.1432
$ws.Range("I132").Value = 3972.7144
$ws.Range("M132").Value = -9388.143199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3222
$ws.Range("K68").Value = 3222
$ws.Range("I68").Value = 3222
$ws.Range("M68").Value = -2473
$ws.Range("H71").Value = 3222
$ws.Range("K71").Value = 16110
$ws.Range("I71").Value = 3222
$ws.Range("M71").Value = -12366
$ws.Range("H132").Value = 5599.6
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 17998.5
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5999.5
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -23058.5
$ws.Range("H136").Value = 6112.909
$ws.Range("K136").Value = 8460.856800000001
$ws.Range("L136").Value = 35625
$ws.Range("I136").Value = 2820.2856
$ws.Range("J136").Value = 11875
$ws.Range("M136").Value = -5910.856800000001
$ws.Range("N136").Value = -40725
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6694499.5
$ws.Range("K2").Value = 6807466
$ws.Range("I2").Value = 6807466
$ws.Range("M2").Value = -6807354
$ws.Range("H28").Value = 16012.5
$ws.Range("L28").Value = 16012.5
$ws.Range("J28").Value = 16012.5
$ws.Range("N28").Value = -16708.5
$ws.Range("H69").Value = 29962.666
$ws.Range("L69").Value = 29962.666
$ws.Range("J69").Value = 29962.666
$ws.Range("N69").Value = -31460.666
$ws.Range("H72").Value = 29962.666
$ws.Range("L72").Value = 89887.99800000001
$ws.Range("J72").Value = 29962.666
$ws.Range("N72").Value = -97375.99800000001
$ws.Range("H137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("N137").Value = -90200
